$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 30727.97
$ws.Range("B3").Value = 1820.37

$ws.Range("A4").Value = "BNB"
$ws.Range("B4").Value = 283.88
$ws.Range("C4").Value = 9
$ws.Range("D4").Value = 2558.7
$ws.Range("E4").Value = "19/07/21"
$ws.Range("F4").Value = "12:58:49"
